$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

$oldHeaders2 = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")
$newHeaders2 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt $oldHeaders2.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders2[$i]
}

# --- 2. Convert the range into an Excel Table (ListObject) ---
# The header row already carries bold/fill/border formatting. If that
# formatting is still present at the moment the table is created, Excel
# captures it into a headerRowDxfId/dxf entry. To keep styles.xml identical
# (dxfs count="0"), stash the header formatting on a scratch cell, blank
# the header formatting, create the table, then restore the formatting via
# a format-only paste (keeps the same style index instead of fragmenting
# the cellXfs table).
$scratch = $ws.Range("W1")
$ws.Range("A1").Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null

$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$rng = $ws.Range("A1:U71")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
